# Apply the "Added CAB checklist & D&A onboarding package" edit to Sheet1.
#
# Summary of the change:
#   1. A new row is inserted at row 214 (pushing the existing rows 214-225
#      down to 215-226) containing a new "D&A Onboarding Package" FAQ entry
#      under Category="D&A Project Mgmt.", Subcategory="Onboarding".
#   2. A brand new row 227 is appended at the bottom of the table with a
#      "D4GV CAB Check List" FAQ entry under Category="D4GV",
#      Subcategory="Project Execution".
#   3. The sheet's AutoFilter range and the workbook-level
#      `_xlnm._FilterDatabase` defined name are both extended from
#      A1:G225 to A1:G226 to account for the extra row.
#   4. The current selection on Sheet1 is moved to E229.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1) Insert the new row 214 ("D&A Onboarding Package") and push the
#    block of rows that used to be 214-225 down to 215-226.
# ---------------------------------------------------------------------
$ws.Rows.Item(214).Insert()

# Match the look of the row directly above (213) which already carries
# the "boxed" style used throughout this table (A/B/C border box, D
# wrapped+bordered, E hyperlink+bordered+wrapped, F bordered).
$ws.Range("A213:F213").Copy()
$ws.Range("A214:F214").PasteSpecial(-4122)  # xlPasteFormats

# Values are written in the same order the original authoring tool
# allocated shared-string ids (D, then F, then E) so the rebuilt
# sharedStrings table lines up with the target workbook.
$ws.Cells.Item(214, 1).Value = "D&A Project Mgmt."
$ws.Cells.Item(214, 2).Value = "Onboarding"
$ws.Cells.Item(214, 3).Value = "General"
$ws.Cells.Item(214, 4).Value = 'Where can I find the document that will provide new suppliers, contractors, and internal new Mondelez employees with information about the setup of the Data and Analytics organization, the multiple teams involved in transforming data into information, as well as the policies, procedures, and standards established within the Mondelez Data and Analytics organization?'
$ws.Cells.Item(214, 6).Value = "D&A Onboarding Package"
$ws.Cells.Item(214, 5).Value = 'https://teams.mdlz.com/:w:/r/sites/dataandanalyticsgrowprogram/Shared%20Documents/Onboarding/D%26A%20Onboarding%20Package.docx?d=wffb7ff4307d44e199e81c2ce25dfc252&csf=1&web=1&e=U81Fw5'

$ws.Rows.Item(214).RowHeight = 87

# ---------------------------------------------------------------------
# 2) Append the new row 227 ("D4GV CAB Check List") at the bottom.
# ---------------------------------------------------------------------
$ws.Range("F216").Copy()
$ws.Range("A227:C227").PasteSpecial(-4122)  # xlPasteFormats -> boxed, non-bold style

$ws.Range("F216").Copy()
$ws.Range("D227").PasteSpecial(-4122)
$ws.Range("D227").Font.Color = 0
$ws.Range("D227").WrapText = $true

$ws.Range("E213").Copy()
$ws.Range("E227").PasteSpecial(-4122)  # hyperlink style, bordered + wrapped

$ws.Range("F12").Copy()
$ws.Range("F227").PasteSpecial(-4122)  # plain hyperlink style

# Same ordering trick as above: A, B, C, D, F, E.
$ws.Cells.Item(227, 1).Value = "D4GV"
$ws.Cells.Item(227, 2).Value = "Project Execution"
$ws.Cells.Item(227, 3).Value = "Project Management"
$ws.Cells.Item(227, 4).Value = "Where can I find the check list for CAB review?"
$ws.Cells.Item(227, 6).Value = "D4GV CAB Check List.xlsx"
$ws.Cells.Item(227, 5).Value = 'https://teams.mdlz.com/:x:/r/sites/ibsdataandanalytics/Shared%20Documents/D%26A%20Calendar/D4GV%20CAB%20Check%20List.xlsx?d=wa29ba62f00ba4fb1921700f2de347e92&csf=1&web=1&e=zXLJvY'

# ---------------------------------------------------------------------
# 3) Grow the AutoFilter range / _FilterDatabase name from G225 to G226.
# ---------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:G226").AutoFilter(1)

# ---------------------------------------------------------------------
# 4) Move the live selection to E229, matching the saved view state.
# ---------------------------------------------------------------------
$ws.Range("E229").Select()

Write-Host "edit applied"
